$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New metric values (same for every data row, columns B..Q)
# Using [double]"..." casts because scientific-notation numeric literals
# (e.g. 1e-05) are not parsed directly by this PowerShell engine.
$values = @(
    [double]"0.9999549410293126",
    [double]"0.9990639202211722",
    [double]"0.9999996702354116",
    [double]"0.9997904501425914",
    [double]"0.999930243544353",
    [double]"4.206057127643839e-05",
    [double]"0.0008737893843822983",
    [double]"2.314895994410189e-07",
    [double]"9.163552103781839e-05",
    [double]"4.59335053186297e-05",
    [double]"0.0004236862227404525",
    [double]"0.006485412190172526",
    [double]"1.0000514959665",
    [double]"0.00676150961376217",
    [double]"110.1527996107309",
    [double]"165.0022117297999"
)

for ($row = 2; $row -le 26; $row++) {
    for ($col = 2; $col -le 17; $col++) {
        $ws.Cells.Item($row, $col).Value = $values[$col - 2]
    }
}
